$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bottom two transaction rows (rows 3 and 4), keeping the
# header row and the first data row. Delete from the bottom up so the
# row indices of the rows above stay valid.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# Update the remaining data row (row 2) with the new transaction:
# "Entertainment" / "douchebag" / 2023-03-22 / 3000.0 / Checkings.
$ws.Cells.Item(2, 1).Value = "Entertainment"
$ws.Cells.Item(2, 2).Value = "douchebag"

# The date and amount look numeric, so force them to stay plain text
# (matching how the original workbook stored every column as a shared
# string) by applying a text number format before assigning the
# values, then dropping the formatting again so no style is left
# behind on the cells.
$dateAndPrice = $ws.Range("C2:D2")
$dateAndPrice.NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "2023-03-22"
$ws.Cells.Item(2, 4).Value = "3000.0"
$dateAndPrice.ClearFormats()

$ws.Cells.Item(2, 5).Value = "Checkings"
